$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (B3=1, C3=2, A3=B3+C3) -----------------------
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("A3").Formula = "=B3+C3"

# --- New conditional formatting rule on A3: highlight when > 10 --------
# xlCellValue = 1, xlGreater = 5 -> "cellIs" / "greaterThan"
$rngA3 = $ws.Range("A3")
$newRule = $rngA3.FormatConditions.Add(1, 5, 10)
# Bring the new rule to the top of the priority list (priority 1), the
# existing "AAAAA" rule on A1:B1 is pushed down to priority 2.
$newRule.SetFirstPriority()
# "Light Red Fill with Dark Red Text" (the same built-in style already
# used by the existing rule): font FF9C0006 on fill FFFFC7CE.
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615

# --- Leave the cursor parked on A3, matching the saved selection -------
$rngA3.Select() | Out-Null
